$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1.05175724944049, 1.061732250817521, 1.058865471869161, 1.070310756288932, 1.056783439289322, 1.064455512831642, 1.06159653861953, 1.073010933640485, 1.058284192987868),
    @(1.053346629947036, 1.063263929734641, 1.060309036370356, 1.071903066491195, 1.058019915742132, 1.065800065445099, 1.062852626451616, 1.074417661830047, 1.059522425379067),
    @(1.054374224834511, 1.064254536129162, 1.061242541667407, 1.072933029893937, 1.058818792020996, 1.06666908702309, 1.063664303075817, 1.075327047972722, 1.060322436153879),
    @(1.054806032897939, 1.064670876511561, 1.061634855790218, 1.073365946019979, 1.059154356991066, 1.067034192431199, 1.064005275694893, 1.07570915465318, 1.060658477664696),
    @(1.054878524308069, 1.064740775623079, 1.061700719541605, 1.07343863004159, 1.059210683434687, 1.067095481867691, 1.064062511651901, 1.075773300561374, 1.060714884098328),
    @(1.054379995422263, 1.064260099715057, 1.061247784296301, 1.072938814847315, 1.058823276959312, 1.066673966479886, 1.06366886016781, 1.075332154475308, 1.060326927461322),
    @(1.052294567142695, 1.062249995217485, 1.059353454885098, 1.070848965684458, 1.057201565920323, 1.064910120493696, 1.062021270095993, 1.073486526459868, 1.058702913406644),
    @(1.048612961053823, 1.058703836005056, 1.056010679151948, 1.067163228654944, 1.054334381240271, 1.061794087246478, 1.059109332154391, 1.070227390372365, 1.055831656995476),
    @(1.046153446737526, 1.056336514650686, 1.053778560551301, 1.064703473515626, 1.052416128232753, 1.059710992844305, 1.057161818270378, 1.06804955802481, 1.053910679848661),
    @(1.045087120018198, 1.055310573094482, 1.052811079498882, 1.063637650253417, 1.05158380340877, 1.058807535004637, 1.056316960302156, 1.067105232086856, 1.053077173027958),
    @(1.04469082716552, 1.05492935224563, 1.052451561657781, 1.06324163671049, 1.051274376567358, 1.058471722419087, 1.056002898748742, 1.066754262157248, 1.05276730676494),
    @(1.04477584306878, 1.05501113191208, 1.052528686479751, 1.063326588507861, 1.051340761760533, 1.058543765812632, 1.0560702771675, 1.06682955582101, 1.052833786232704),
    @(1.04505436670378, 1.055279064177524, 1.052781364796674, 1.063604918138077, 1.051558231526437, 1.058779781332343, 1.056291004891294, 1.067076225038379, 1.05305156483062),
    @(1.045225945958517, 1.055444127207644, 1.052937027818698, 1.063776390121101, 1.051692186569135, 1.05892516781143, 1.056426970008623, 1.067228178637249, 1.053185710104843),
    @(1.046224186212693, 1.056404583803001, 1.053842748177146, 1.064774192467354, 1.052471330240454, 1.059770920823153, 1.057217854954322, 1.068112201508568, 1.053965960249545),
    @(1.046849989808694, 1.057006812114233, 1.054410620095457, 1.06539988478618, 1.052959603656946, 1.060301041604038, 1.057713530093708, 1.068666369025146, 1.054454927070293),
    @(1.047214881783904, 1.057357997803006, 1.054741758210292, 1.065764770425395, 1.053244240942755, 1.06061011136363, 1.05800249789374, 1.068989479720653, 1.054739968573697),
    @(1.047339278999492, 1.057477729085189, 1.054854652398585, 1.065889175384912, 1.053341267257324, 1.060715472540578, 1.058101003040974, 1.069099631025886, 1.054837132676762),
    @(1.046782860388658, 1.056942207458152, 1.054349702385181, 1.065332761239746, 1.052907233593887, 1.060244179275871, 1.05766036452848, 1.068606925172389, 1.054402482635738),
    @(1.044972354359532, 1.055200168752938, 1.052706961600791, 1.063522960367338, 1.051494199459454, 1.058710287007088, 1.056226012866135, 1.06700359283498, 1.052987441830757),
    @(1.043832787831138, 1.054104062682096, 1.051673219985604, 1.062384371735762, 1.050604235588034, 1.057744544170503, 1.055322765362694, 1.065994322560932, 1.052096214108559),
    @(1.044437013574077, 1.054685209695664, 1.05222131265032, 1.062988028074642, 1.051076170171822, 1.058256630939724, 1.055801730495723, 1.066529471595368, 1.052568818893588),
    @(1.046813193681789, 1.056971399780816, 1.054377228760683, 1.065363091697928, 1.052930897879111, 1.06026987334938, 1.057684388226073, 1.068633785679766, 1.054426180526962),
    @(1.049565606285847, 1.059621135994137, 1.056875470313181, 1.068116502766107, 1.055076786112191, 1.062600637123987, 1.059863205819263, 1.071070817747025, 1.056575116167469),
)

$cols = @("C", "D", "E", "F", "J", "K", "L", "M", "N")
$startRow = 2

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range($cols[$j] + $row).Value = $data[$i][$j]
    }
}